$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 1052.6666
$ws.Range("I4").Value = 920.5714
$ws.Range("J4").Value = 1515
$ws.Range("K4").Value = 920.5714
$ws.Range("L4").Value = 1515
$ws.Range("M4").Value = -806.5714
$ws.Range("N4").Value = -1743
$ws.Range("H19").Value = 2215.077
$ws.Range("I19").Value = 2644.6
$ws.Range("J19").Value = 783.3333
$ws.Range("K19").Value = 2644.6
$ws.Range("L19").Value = 783.3333
$ws.Range("M19").Value = -2469.6
$ws.Range("N19").Value = -1133.3333
$ws.Range("H62").Value = 3372.5
$ws.Range("I62").Value = 3650
$ws.Range("J62").Value = 3095
$ws.Range("K62").Value = 3650
$ws.Range("L62").Value = 3095
$ws.Range("M62").Value = -3026
$ws.Range("N62").Value = -4343
$ws.Range("H65").Value = 3372.5
$ws.Range("I65").Value = 3650
$ws.Range("J65").Value = 3095
$ws.Range("K65").Value = 18250
$ws.Range("L65").Value = 15475
$ws.Range("M65").Value = -15130
$ws.Range("N65").Value = -21715
$ws.Range("H112").Value = 2084.6667
$ws.Range("J112").Value = 2171.48
$ws.Range("L112").Value = 6514.440000000001
$ws.Range("N112").Value = -8730.440000000001
$ws.Range("H132").Value = 2840.9697
$ws.Range("I132").Value = 2733.9355
$ws.Range("J132").Value = 4500
$ws.Range("K132").Value = 8201.806500000001
$ws.Range("L132").Value = 13500
$ws.Range("M132").Value = -5671.806500000001
$ws.Range("N132").Value = -18560
$ws.Range("H135").Value = 4699.6
$ws.Range("I135").Value = 5146
$ws.Range("K135").Value = 46314
$ws.Range("M135").Value = -43779
$ws.Range("H137").Value = 23317.732
$ws.Range("I137").Value = 15998.75
$ws.Range("K137").Value = 47996.25
$ws.Range("M137").Value = -45446.25
$ws.Range("H138").Value = 5423.3335
$ws.Range("J138").Value = 5581.9062
$ws.Range("L138").Value = 16745.7186
$ws.Range("N138").Value = -27025.7186

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3793
$ws.Range("I2").Value = 3404.718
$ws.Range("J2").Value = 4957.846
$ws.Range("K2").Value = 3404.718
$ws.Range("L2").Value = 4957.846
$ws.Range("M2").Value = -3291.718
$ws.Range("N2").Value = -5183.846
$ws.Range("H32").Value = 12099.742
$ws.Range("I32").Value = 10622.339
$ws.Range("K32").Value = 10622.339
$ws.Range("M32").Value = -10335.339
$ws.Range("H45").Value = 6556
$ws.Range("I45").Value = 7310.684
$ws.Range("K45").Value = 7310.684
$ws.Range("M45").Value = -6933.684
$ws.Range("H61").Value = 7909.6816
$ws.Range("I61").Value = 3369.2104
$ws.Range("J61").Value = 36666
$ws.Range("K61").Value = 3369.2104
$ws.Range("L61").Value = 36666
$ws.Range("M61").Value = -3157.2104
$ws.Range("N61").Value = -37090
$ws.Range("H74").Value = 14090.086
$ws.Range("I74").Value = 14359.08
$ws.Range("J74").Value = 13417.6
$ws.Range("K74").Value = 14359.08
$ws.Range("L74").Value = 13417.6
$ws.Range("M74").Value = -13485.08
$ws.Range("N74").Value = -15165.6
$ws.Range("H77").Value = 14090.086
$ws.Range("I77").Value = 14359.08
$ws.Range("J77").Value = 13417.6
$ws.Range("K77").Value = 71795.39999999999
$ws.Range("L77").Value = 67088
$ws.Range("M77").Value = -67427.39999999999
$ws.Range("N77").Value = -75824
$ws.Range("H88").Value = 2400
$ws.Range("J88").Value = 2457.1428
$ws.Range("L88").Value = 2457.1428
$ws.Range("N88").Value = -3269.1428
$ws.Range("H91").Value = 2400
$ws.Range("J91").Value = 2457.1428
$ws.Range("L91").Value = 2457.1428
$ws.Range("N91").Value = -5265.1428
$ws.Range("H97").Value = 1293.069
$ws.Range("I97").Value = 1202.68
$ws.Range("J97").Value = 1858
$ws.Range("K97").Value = 1202.68
$ws.Range("L97").Value = 1858
$ws.Range("M97").Value = -706.6800000000001
$ws.Range("N97").Value = -2850
$ws.Range("H102").Value = 5595.5557
$ws.Range("I102").Value = 1295.125
$ws.Range("J102").Value = 39999
$ws.Range("K102").Value = 1295.125
$ws.Range("L102").Value = 39999
$ws.Range("M102").Value = 326.875
$ws.Range("N102").Value = -43243
$ws.Range("H110").Value = 1786.5
$ws.Range("I110").Value = 1679.2106
$ws.Range("K110").Value = 1679.2106
$ws.Range("M110").Value = 365.7893999999999
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").Value = ""
$ws.Range("H116").Value = 3793
$ws.Range("I116").Value = 3404.718
$ws.Range("J116").Value = 4957.846
$ws.Range("K116").Value = 3404.718
$ws.Range("L116").Value = 4957.846
$ws.Range("M116").Value = -1110.718
$ws.Range("N116").Value = -9545.846
$ws.Range("H122").Value = 3537.2
$ws.Range("I122").Value = 3114.2307
$ws.Range("K122").Value = 9342.6921
$ws.Range("M122").Value = -6892.6921
$ws.Range("H132").Value = 2578.7058
$ws.Range("J132").Value = 2200
$ws.Range("L132").Value = 6600
$ws.Range("N132").Value = -11660
$ws.Range("H136").Value = 7909.6816
$ws.Range("I136").Value = 3369.2104
$ws.Range("J136").Value = 36666
$ws.Range("K136").Value = 10107.6312
$ws.Range("L136").Value = 109998
$ws.Range("M136").Value = -7557.6312
$ws.Range("N136").Value = -115098

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3793
$ws.Range("I3").Value = 3404.718
$ws.Range("J3").Value = 4957.846
$ws.Range("K3").Value = 3404.718
$ws.Range("L3").Value = 4957.846
$ws.Range("M3").Value = -3290.718
$ws.Range("N3").Value = -5185.846
$ws.Range("H86").Value = 669875.0600000001
$ws.Range("I86").Value = 1251646.6
$ws.Range("J86").Value = 4993.2856
$ws.Range("K86").Value = 1251646.6
$ws.Range("L86").Value = 4993.2856
$ws.Range("M86").Value = -1250523.6
$ws.Range("N86").Value = -7239.2856
$ws.Range("H89").Value = 669875.0600000001
$ws.Range("I89").Value = 1251646.6
$ws.Range("J89").Value = 4993.2856
$ws.Range("K89").Value = 6258233
$ws.Range("L89").Value = 24966.428
$ws.Range("M89").Value = -6252617
$ws.Range("N89").Value = -36198.428
$ws.Range("H94").Value = 1553.3334
$ws.Range("I94").Value = 1622.7142
$ws.Range("J94").Value = 1067.6666
$ws.Range("K94").Value = 1622.7142
$ws.Range("L94").Value = 1067.6666
$ws.Range("M94").Value = -1171.7142
$ws.Range("N94").Value = -1969.6666
$ws.Range("H134").Value = 11519.03
$ws.Range("I134").Value = 5619.2173
$ws.Range("K134").Value = 16857.6519
$ws.Range("M134").Value = -14322.6519

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 866.8929000000001
$ws.Range("I22").Value = 727.8261
$ws.Range("K22").Value = 727.8261
$ws.Range("M22").Value = -377.8261
$ws.Range("H31").Value = 5642.4287
$ws.Range("I31").Value = 4000
$ws.Range("J31").Value = 5916.1665
$ws.Range("K31").Value = 4000
$ws.Range("L31").Value = 5916.1665
$ws.Range("M31").Value = -3705
$ws.Range("N31").Value = -6506.1665
$ws.Range("H34").Value = 5642.4287
$ws.Range("I34").Value = 4000
$ws.Range("J34").Value = 5916.1665
$ws.Range("K34").Value = 4000
$ws.Range("L34").Value = 5916.1665
$ws.Range("M34").Value = -3798
$ws.Range("N34").Value = -6320.1665
$ws.Range("H58").Value = 5905.5435
$ws.Range("I58").Value = 4689.154
$ws.Range("J58").Value = 7486.85
$ws.Range("K58").Value = 4689.154
$ws.Range("L58").Value = 7486.85
$ws.Range("M58").Value = -4486.154
$ws.Range("N58").Value = -7892.85
$ws.Range("H107").Value = 679.2121
$ws.Range("I107").Value = 707.6539
$ws.Range("J107").Value = 573.5714
$ws.Range("K107").Value = 707.6539
$ws.Range("L107").Value = 573.5714
$ws.Range("M107").Value = 1212.3461
$ws.Range("N107").Value = -4413.5714
$ws.Range("H122").Value = 5702.615
$ws.Range("I122").Value = 5925.2
$ws.Range("J122").Value = 5563.5
$ws.Range("K122").Value = 17775.6
$ws.Range("L122").Value = 16690.5
$ws.Range("M122").Value = -15325.6
$ws.Range("N122").Value = -21590.5
$ws.Range("H136").Value = 5905.5435
$ws.Range("I136").Value = 4689.154
$ws.Range("J136").Value = 7486.85
$ws.Range("K136").Value = 14067.462
$ws.Range("L136").Value = 22460.55
$ws.Range("M136").Value = -11517.462
$ws.Range("N136").Value = -27560.55

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 1900.65
$ws.Range("J12").Value = 2527.5334
$ws.Range("L12").Value = 7582.600199999999
$ws.Range("N12").Value = -7928.600199999999
$ws.Range("H137").Value = 6321.4736
$ws.Range("I137").Value = 6545.5
$ws.Range("J137").Value = 5937.4287
$ws.Range("K137").Value = 19636.5
$ws.Range("L137").Value = 17812.2861
$ws.Range("M137").Value = -14536.5
$ws.Range("N137").Value = -28012.2861

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2217.1082
$ws.Range("I102").Value = 2266.0293
$ws.Range("K102").Value = 2266.0293
$ws.Range("M102").Value = -644.0293000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1697.8667
$ws.Range("I46").Value = 1086.5
$ws.Range("J46").Value = 1920.1818
$ws.Range("K46").Value = 1086.5
$ws.Range("L46").Value = 1920.1818
$ws.Range("M46").Value = -898.5
$ws.Range("N46").Value = -2296.1818
$ws.Range("H61").Value = 2203.25
$ws.Range("I61").Value = 1758.909
$ws.Range("K61").Value = 1758.909
$ws.Range("M61").Value = -1556.909
$ws.Range("H113").Value = 2203.25
$ws.Range("I113").Value = 1758.909
$ws.Range("K113").Value = 1758.909
$ws.Range("M113").Value = 411.0909999999999
$ws.Range("H136").Value = 6561.483
$ws.Range("I136").Value = 6187.1094
$ws.Range("J136").Value = 8269.5625
$ws.Range("K136").Value = 18561.3282
$ws.Range("L136").Value = 24808.6875
$ws.Range("M136").Value = -16011.3282
$ws.Range("N136").Value = -29908.6875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 16263
$ws.Range("J14").Value = 3000
$ws.Range("L14").Value = 3000
$ws.Range("N14").Value = -3336
$ws.Range("H113").Value = 649.1852
$ws.Range("I113").Value = 602.1
$ws.Range("K113").Value = 1806.3
$ws.Range("M113").Value = 363.6999999999998
$ws.Range("H126").Value = 5781.7617
$ws.Range("I126").Value = 3678.5518
$ws.Range("K126").Value = 11035.6554
$ws.Range("M126").Value = -8565.6554

